# Update the localization-status workbook to mark the e6ead3b7 file as
# "Ready for handoff" (instead of "In Translation") and refresh the
# handoff timestamps / priority for that row.

$wb = $excel.ActiveWorkbook

$newStatus   = "Ready for handoff"
$newPriority = "mt"
$zhDateTime  = "2016-08-22 11:02:57"
$deDateTime  = "2016-08-22 11:03:05"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $deDateTime

# --- zh-cn sheet ------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("E3").Value = $newPriority
$wsZh.Range("H3").Value = $zhDateTime

# --- de-de sheet ------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("E3").Value = $newPriority
$wsDe.Range("H3").Value = $deDateTime

# --- Resize the status columns to fit the longer "Ready for handoff" text --
$newWidth = 17.2159881591797
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZh.Columns.Item(3).ColumnWidth = $newWidth
$wsDe.Columns.Item(3).ColumnWidth = $newWidth
